{"js": "// Indent every wrapped line (the text following each manual line break\n// `<w:br/>`) in the \"Vitals\" paragraph and in the \"Summary Statement\"\n// paragraph by prepending four spaces. The first line of each paragraph is\n// left untouched.\n//\n// Approach: locate each target paragraph by a distinctive snippet of its\n// first line, read its full text (Office.js represents a manual line break\n// as \"\\v\"), rebuild the line list with four spaces prepended to every line\n// but the first, then push the new content back as OOXML (<w:r>/<w:t>/\n// <w:br/> run) via insertOoxml(\"Replace\") on the paragraph's range. Using\n// raw OOXML (rather than Range.insertText) avoids the shim's side effect of\n// re-marking every <w:t> in the run with xml:space=\"preserve\" when only one\n// segment actually changed.\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\n// Build the replacement <w:p> children (run content only -- no <w:pPr>,\n// since the paragraphs we touch here carry none) from a paragraph's current\n// text, indenting every line after the first by `indent` spaces.\nfunction buildIndentedRunOoxml(paragraphText, indent) {\n  const lines = paragraphText.split(\"\\v\");\n  const pieces = lines.map((line, i) => {\n    const text = i === 0 ? line : indent + line;\n    // Preserve leading/trailing whitespace explicitly, matching Word's own\n    // behavior of adding xml:space=\"preserve\" whenever it matters.\n    const needsPreserve = /^\\s|\\s$/.test(text) || text === \"\";\n    const openTag = needsPreserve ? '<w:t xml:space=\"preserve\">' : \"<w:t>\";\n    return openTag + xmlEscape(text) + \"</w:t>\";\n  });\n  return \"<w:r>\" + pieces.join(\"<w:br/>\") + \"</w:r>\";\n}\n\nfunction wrapParagraphPackageXml(paragraphInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" +\n    paragraphInnerXml +\n    \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nasync function indentWrappedLines(context, anchorText, indent) {\n  const body = context.document.body;\n  const results = body.search(anchorText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not locate paragraph containing: \" + anchorText);\n  }\n\n  const para = results.items[0].paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n\n  const newRunXml = buildIndentedRunOoxml(para.text, indent);\n  const pkgXml = wrapParagraphPackageXml(newRunXml);\n\n  const range = para.getRange();\n  range.insertOoxml(pkgXml, \"Replace\");\n  await context.sync();\n}\n\n// \"Vitals\" paragraph: Heart Rate / Respiratory Rate / Weight lines.\nawait indentWrappedLines(context, \"Heart Rate: 101, Blood Pressure: 197/86\", \"    \");\n\n// \"Summary Statement\" paragraph.\nawait indentWrappedLines(\n  context,\n  \"This is a 49 year old female, who is presenting today for RLQ pain radiating down left leg for past 12 hours.\",\n  \"    \"\n);\n", "ps1": "# Indent every wrapped line (the text following each manual line break) in\n# the \"Vitals\" paragraph and in the \"Summary Statement\" paragraph by\n# prepending four spaces. The first line of each paragraph is left\n# untouched.\n#\n# Approach: locate each target paragraph by a distinctive snippet at the\n# start of its text, split its Range.Text on Chr(11) (Word's manual-line-\n# break character, cf. <w:br/>), rebuild the line list with four spaces\n# prepended to every line but the first, then push the new content back as\n# WordOpenXML via Range.InsertXML -- this replaces only the run content of\n# that paragraph without touching unrelated runs, avoiding Word re-marking\n# every <w:t> in the run with xml:space=\"preserve\" the way a plain\n# Range.Text/.InsertBefore edit would.\n\n$d = $word.ActiveDocument\n$vt = [char]11\n$cr = [char]13\n\nfunction Escape-XmlText($text) {\n    $t = $text -replace '&', '&amp;'\n    $t = $t -replace '<', '&lt;'\n    $t = $t -replace '>', '&gt;'\n    $t = $t -replace '\"', '&quot;'\n    return $t\n}\n\nfunction Needs-Preserve($text) {\n    if ($text -eq '') { return $true }\n    if ($text -match '^\\s') { return $true }\n    if ($text -match '\\s$') { return $true }\n    return $false\n}\n\nfunction Find-ParagraphStartingWith($doc, $prefix) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.StartsWith($prefix)) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Indent-WrappedLines($doc, $prefix, $indent) {\n    $para = Find-ParagraphStartingWith $doc $prefix\n    if ($para -eq $null) {\n        throw \"Could not locate paragraph starting with: $prefix\"\n    }\n\n    $rng = $para.Range\n\n    # Range.Text includes the trailing paragraph mark (Chr(13)); strip it.\n    $full = $rng.Text\n    if ($full.Length -gt 0 -and $full[$full.Length - 1] -eq $cr) {\n        $full = $full.Substring(0, $full.Length - 1)\n    }\n\n    $lines = $full.Split($vt)\n\n    # Preserve this paragraph's own formatting (<w:pPr>), if any, by pulling\n    # it out of the paragraph's own WordOpenXML. WordOpenXML returns the\n    # *whole* mini-package (styles.xml and friends included, which are full\n    # of unrelated <w:pPr> blocks), so this must be scoped down to: (1) the\n    # \"/word/document.xml\" part, (2) its <w:body>, and (3) not crossing past\n    # the first paragraph's own closing </w:p> tag.\n    $pprXml = ''\n    $oxml = $rng.WordOpenXML\n    if ($oxml -match '(?s)pkg:name=\"/word/document\\.xml\"[^>]*>.*?<pkg:xmlData>(?<doc>.*?)</pkg:xmlData>') {\n        $docXml = $matches['doc']\n        if ($docXml -match '(?s)<w:body>(?<body>.*?)</w:body>') {\n            $bodyXml = $matches['body']\n            if ($bodyXml -match '(?s)^<w:p[ >](?:(?!</w:p>).)*?(?<ppr><w:pPr>.*?</w:pPr>)') {\n                $pprXml = $matches['ppr']\n            }\n        }\n    }\n\n    $runXml = '<w:r>'\n    for ($i = 0; $i -lt $lines.Length; $i++) {\n        if ($i -gt 0) {\n            $runXml += '<w:br/>'\n        }\n        $line = $lines[$i]\n        if ($i -gt 0) {\n            $line = $indent + $line\n        }\n        $escaped = Escape-XmlText $line\n        if (Needs-Preserve $line) {\n            $runXml += '<w:t xml:space=\"preserve\">' + $escaped + '</w:t>'\n        } else {\n            $runXml += '<w:t>' + $escaped + '</w:t>'\n        }\n    }\n    $runXml += '</w:r>'\n\n    $pkgXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' + $pprXml + $runXml + '</w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $rng.InsertXML($pkgXml)\n}\n\n# \"Vitals\" paragraph: Heart Rate / Respiratory Rate / Weight lines.\nIndent-WrappedLines $d \"Heart Rate: 101, Blood Pressure: 197/86\" \"    \"\n\n# \"Summary Statement\" paragraph.\nIndent-WrappedLines $d \"This is a 49 year old female, who is presenting today for RLQ pain radiating down left leg for past 12 hours.\" \"    \"\n"}
